$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.102614641189575
$ws.Range("B1").Value = 1.639281988143921
$ws.Range("C1").Value = 4.62448787689209
$ws.Range("D1").Value = 0.4045725464820862
$ws.Range("E1").Value = 0.4551762342453003
